$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.04763786555579896
$ws.Range("C2").Value = 114.8270160096505
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 645.3272768299601
$ws.Range("G2").Value = 761.0074203417504

# Row 3
$ws.Range("B3").Value = 3.230985683306322
$ws.Range("C3").Value = 114.8270160096505
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 645.3272768299601
$ws.Range("G3").Value = 764.1907681595009
